$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRow = 84

$ws.Cells.Item($newRow, 1).Value = "2025-06"
$ws.Cells.Item($newRow, 2).Value = 1
$ws.Cells.Item($newRow, 3).Value = 252
$ws.Cells.Item($newRow, 4).Value = 0.3968253968253968
